$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: update label and count, keep existing checkbox/boolean value
$ws.Range("A1").Value = "ar2"
$ws.Range("C1").Value = 30

# Row 2: turn into "Lâmpada" entry with new count and boolean flag column D (E2 removed)
$ws.Range("A2").Value = "a"
$ws.Range("B2").Value = "Lâmpada"
$ws.Range("C2").Value = 100
$ws.Range("D2").Value = $false
$ws.Range("E2").ClearContents()

# Row 3: turn into "A/C" entry with new count
$ws.Range("A3").Value = "ae"
$ws.Range("B3").Value = "A/C"
$ws.Range("C3").Value = 23
$ws.Range("D3").Value = $false

# Row 4 (new row): "Lâmpada" entry
$ws.Range("A4").Value = "la"
$ws.Range("B4").Value = "Lâmpada"
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = $false
